$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data set for rows 2..19 (columns A:D)
# Column A: Residue Combination, B: Count, C: Permeated Residues, D: Frames
$data = @(
    @("130, 780",                       1, "130",       "3631"),
    @("455",                            1, "455",       "3666"),
    @("1073, 1105",                     1, "1105",      "5131"),
    @("130, 1073",                      1, "130",       "4415"),
    @("780, 1073, 1105",                1, "1105",      "4994"),
    @("130, 455, 780",                  2, "130, 130",  "5269, 6424"),
    @("130, 1073, 1105",                1, "1105",      "5399"),
    @("423, 748, 780, 1073",            2, "780, 780",  "5677, 5887"),
    @("423, 1105, 1105",                1, "1105",      "5331"),
    @("423, 1073, 1105",                1, "1105",      "5433"),
    @("98, 130, 748, 780, 1073",        1, "130",       "6016"),
    @("98, 130, 748, 1073",             2, "130, 130",  "5582, 6488"),
    @("423, 748, 780, 1073, SF",        1, "780",       "6202"),
    @("130, 423, 748, 1073",            1, "130",       "6561"),
    @("130, 780, 1073, 1105",           1, "1105",      "6359"),
    @("98, 130, 130, 423, 780",         1, "130",       "6727"),
    @("98, 130, 423, 1073",             1, "130",       "6670"),
    @("98, 130, 455, 780",              1, "780",       "6748")
)

$lastRow = 1 + $data.Count

# Columns C and D always hold text (residue / frame lists), even when a
# given row's text happens to look like a bare number (e.g. "455", "3631").
# Force text formatting first so Excel stores these as strings rather than
# silently coercing them to numeric cells, matching the source data.
$ws.Range("C2:D$lastRow").NumberFormat = "@"

# Column A is normally a comma-separated residue list (non-numeric text),
# except for the single-residue row "455" (row 3), which Excel would
# otherwise auto-coerce to a number. Force text format just for that cell
# so it round-trips as a string like the rest of the column.
$ws.Cells.Item(3, 1).NumberFormat = "@"

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $rowIndex++
}
